# "Improved log in the TA catch"
# The SHHConfig row previously logged the literal text "ssh" for every
# column; change it to point at the actual SSH config file so the log
# line is actually useful.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 holds the "SHHConfig" entries (B7:AW7) -- replace "ssh" with the
# path to the SSH config file used by the catch/log.
$ws.Range("B7:AW7").Value = "Configurations/SSHConfig.csv"

# Reflect the edit in the sheet view: select the row that was just
# changed and scroll it into view (matches the saved view state of the
# authored workbook).
$ws.Activate()
$ws.Range("B7:AW7").Select()
$excel.ActiveWindow.ScrollColumn = 49
